# Auto-save via app Streamlit
# A new reservation row ("Fermeture" / "Autre") is inserted above the old
# row 32 (Mika Thielen), pushing every following row down by one. The
# Luna Ciccardi phone number also loses its leading "+" and becomes a
# plain number instead of text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32; everything from the old row 32 down
# (including the TOTAL row) shifts down by one row.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the "Fermeture" entry.
$ws.Range("A32").Value = "Fermeture"
$ws.Range("B32").Value = "Autre"
$ws.Range("D32").Value = 45878
$ws.Range("E32").Value = 45880
$ws.Range("F32").Value = 2
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2025
$ws.Range("L32").Value = 8

# The date-insert carried the neighbouring row's "last-modified" date
# style (s="3") into N32/O32; the new row has no such timestamps, so
# drop the inherited formatting.
$ws.Range("N32").ClearFormats()
$ws.Range("O32").ClearFormats()

# Luna Ciccardi's row (now row 35 after the shift) gets her phone number
# normalised from text "+393485814122" to the plain number 393485814122.
$ws.Range("C35").Value = 393485814122
